# Insert a new data row at row 152 (pushing existing rows 152..207 down to
# 153..208), then populate the new row with its own values, matching the
# weekly refresh of the "Hortaliza, Vega Modelo de Temuco - Espinaca" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 152:207 down by one row, creating a blank row 152.
$ws.Rows("152:152").Insert()

# Populate the newly inserted row 152 with the new record.
$ws.Cells.Item(152, 1).Value = 10
$ws.Cells.Item(152, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(152, 3).Value = "La Araucanía"
$ws.Cells.Item(152, 4).Value = 44875
$ws.Cells.Item(152, 5).Value = 9
$ws.Cells.Item(152, 6).Value = 100112012
$ws.Cells.Item(152, 7).Value = "Espinaca"
$ws.Cells.Item(152, 8).Value = "Sin especificar"
$ws.Cells.Item(152, 9).Value = "Primera"
$ws.Cells.Item(152, 10).Value = 110
$ws.Cells.Item(152, 11).Value = 9000
$ws.Cells.Item(152, 12).Value = 9000
$ws.Cells.Item(152, 13).Value = 9000
$ws.Cells.Item(152, 14).Value = "$/docena de atados"
$ws.Cells.Item(152, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(152, 16).Value = 3000
$ws.Cells.Item(152, 17).Value = 3
$ws.Cells.Item(152, 18).Value = "Hortaliza"
